$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.511.37"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "1.976.57"
$ws.Range("E3").Value = "  +3.90%  "
$ws.Range("E4").Value = "  +0.42%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "327.14"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.003"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.26%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4661"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("E8").Value = "  -0.36%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "46.15"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.44%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07950"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.80%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.9930"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.37%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "22.88"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +3.91%  "
$ws.Range("D13").Value = "1.980.53"
$ws.Range("E13").Value = "  +4.15%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.195"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +1.61%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.852"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.71%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.07088"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "87.74"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000009936"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.52%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.30"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +0.37%  "
$ws.Range("D22").Value = "29.518.57"
$ws.Range("E22").Value = "  +0.75%  "
$ws.Range("B23").Value = "BitDAO"
$ws.Range("C23").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.5062"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +5.29%  "
$ws.Range("B24").Value = "Uniswap"
$ws.Range("C24").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.554"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("B25").Value = "Cosmos"
$ws.Range("C25").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "11.18"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.68%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.217.97"
$ws.Range("E26").Value = "  +4.27%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "2.108"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "158.32"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "19.56"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.73%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.798"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.25%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "119.59"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("B32").Value = "LidoDAOToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.910"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.09425"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.8944"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "5.236"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.323"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.189"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.05817"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.170"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -1.36%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.02104"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "7.779"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.5720"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.11%  "
$ws.Range("B43").Value = "PEPE"
$ws.Range("C43").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.000003149"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +43.92%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.1801"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.79%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "9.660"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.765"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.69%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "11.83"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("B48").Value = "Decentraland"
$ws.Range("C48").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.5364"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.23%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "2.187"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06926"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.63%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "114.10"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.46%  "
